# Apply "Updated symbol list on Mon Dec 26 22:32:53 UTC 2022 with GitHub Actions"
# This script updates the Price column (D) for most rows, and for rows 10-18
# it shifts the Coin/Link/Volume(1h) data up by one (new coin appended at the
# bottom of that block) along with new Price values.
#
# Price cells are stored as plain text in the workbook (not numbers), so we
# force NumberFormat to "@" (Text) before assigning each one. Otherwise Excel
# auto-detects the numeric-looking string and stores it as a real number,
# which both changes the cell type and can mangle the literal text
# (trailing zeros / leading zeros / scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Simple Price (column D) updates ---
Set-TextValue "D2"  "242.73"
Set-TextValue "D3"  "23.09"
Set-TextValue "D4"  "5.416"
Set-TextValue "D5"  "0.05903"
Set-TextValue "D7"  "6.531"
Set-TextValue "D8"  "0.8096"
Set-TextValue "D9"  "0.9365"

Set-TextValue "D19" "0.006012"
Set-TextValue "D21" "0.004898"
Set-TextValue "D22" "0.00006806"
Set-TextValue "D23" "3.570"
Set-TextValue "D24" "2.142"

Set-TextValue "D40" "0.03960"
Set-TextValue "D41" "0.006476"
Set-TextValue "D42" "0.1074"
Set-TextValue "D43" "0.002572"
Set-TextValue "D44" "0.008789"
Set-TextValue "D45" "0.00005235"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.6707"
Set-TextValue "D48" "0.002392"
Set-TextValue "D49" "0.00002102"
Set-TextValue "D50" "0.0002002"

# --- Rows 10-18: Coin/Link/Volume(1h) shifted up by one row, new Price values ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1427"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07418"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03264"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03089"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09357"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.873"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001596"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04687"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005906"
$ws.Range("E18").Value = "17OneONE"
